$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated demand values for column B, rows 2-25 (hour index 0-23 in column A)
$newValues = @{
    2  = 202
    3  = 450
    4  = 957
    5  = 3186
    6  = 14096
    7  = 49841
    8  = 133857
    9  = 265188
    10 = 208692
    11 = 69506
    12 = 31592
    13 = 23921
    14 = 102380
    15 = 105904
    16 = 123157
    17 = 178291
    18 = 226572
    19 = 199901
    20 = 128416
    21 = 89067
    22 = 67722
    23 = 52021
    24 = 33900
    25 = 19588
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
